$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) that mirrors row 2, but with column B ("属性")
# set to "风" instead of "水".
$ws.Range("A3").Value = 114514
$ws.Range("B3").Value = "风"
$ws.Range("C3").Value = "test"
$ws.Range("D3").Value = "传奇道具"
$ws.Range("E3").Value = "114514水"
$ws.Range("F3").Value = "1水"
$ws.Range("G3").Value = "休息休息"
$ws.Range("H3").Value = "休息休息休息"
$ws.Range("I3").Value = 114514

# Match the final selection left behind in the saved file.
$ws.Range("C8").Select()
